# Apply "Final commit with lot of minor errors solved" edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the misspelled name in B3 (Prretika Shetty -> Preetika Shetty)
$ws.Range("B3").Value = "Preetika Shetty"

# Update column A (employee/record ids): 1..10 -> 52501..52510
$ws.Range("A1").Value = 52501
$ws.Range("A2").Value = 52502
$ws.Range("A3").Value = 52503
$ws.Range("A4").Value = 52504
$ws.Range("A5").Value = 52505
$ws.Range("A6").Value = 52506
$ws.Range("A7").Value = 52507
$ws.Range("A8").Value = 52508
$ws.Range("A9").Value = 52509
$ws.Range("A10").Value = 52510

# Update column C (scores)
$ws.Range("C1").Value = 80
$ws.Range("C2").Value = 80
$ws.Range("C3").Value = 97
$ws.Range("C4").Value = 96
$ws.Range("C5").Value = 93
$ws.Range("C6").Value = 95
$ws.Range("C7").Value = 92
$ws.Range("C8").Value = 85
$ws.Range("C9").Value = 91
$ws.Range("C10").Value = 90

# Update the active selection to B3 (single active cell instead of full range)
$ws.Range("B3").Select()
